$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.56315883341524242,
    0.56447796499544467,
    0.5657613997973141,
    0.56701037894739581,
    0.5682260982607602,
    0.56940970994965845,
    0.5705623239776253,
    0.57168500926470667,
    0.57277879484228877,
    0.57384467100238046,
    0.57488359045896442,
    0.57589646952526719,
    0.57688418930409657,
    0.5778475968854232,
    0.57878750654441502,
    0.57970470093319038,
    0.58059993226013029,
    0.58147392345138216,
    0.5823273692900498,
    0.58316093752939746,
    0.58397526997717242,
    0.58477098354883805,
    0.58554867128811794,
    0.58630890335376518,
    0.58705222797191314,
    0.58777917235373067,
    0.58849024357841018,
    0.58918592944176662,
    0.58986669927092417,
    0.59053300470573467,
    0.59118528044768903,
    0.59182394497719026,
    0.59244940124011836,
    0.59306203730467555,
    0.59366222698953175,
    0.59425033046430964,
    0.59482669482345807,
    0.5953916546345589,
    0.59594553246210413,
    0.59648863936776619,
    0.59702127538815919,
    0.59754372999107064,
    0.59805628251110921,
    0.59855920256569228,
    0.59905275045225825,
    0.5995371775275643,
    0.60001272656989169,
    0.60047963212495248,
    0.60093812083625853,
    0.60138841176067803,
    0.60183071666987975,
    0.60226524033832829,
    0.60269218081846831,
    0.60311172970370186,
    0.60352407237974237,
    0.60392938826489173,
    0.60432785103976949,
    0.60471962886699382,
    0.60510488460129019,
    0.60548377599048187,
    0.60585645586779246,
    0.60622307233587269,
    0.60658376894293942,
    0.60693868485139868,
    0.60728795499930621,
    0.60763171025499962,
    0.60797007756522181,
    0.60830318009703754,
    0.60863113737383279,
    0.60895406540566832,
    0.60927207681425055,
    0.60958528095276487,
    0.60989378402080685,
    0.61019768917463668,
    0.61049709663296714,
    0.61079210377848903,
    0.6110828052553251,
    0.61136929306259646,
    0.61165165664427501,
    0.61192998297548595,
    0.6122043566454215,
    0.61247485993701067,
    0.61274157290349263,
    0.61300457344202564,
    0.61326393736446327,
    0.6135197384654193,
    0.61377204858774015,
    0.61402093768549415,
    0.61426647388458755,
    0.61450872354110408
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

Write-Output "Updated $($values.Length) Lambda values in column A (rows 2-91)"
